$wb = $excel.ActiveWorkbook

# --- Section_A sheet updates ---
$wsA = $wb.Worksheets.Item("Section_A")
$wsA.Range("B2").Value = "CS309"
$wsA.Range("F2").Value = "CS309"
$wsA.Range("C3").Value = "CS304"
$wsA.Range("D3").Value = "CS304"
$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "CS309"
$wsA.Range("E5").Value = "CS304"
$wsA.Range("F5").Value = "Free"
$wsA.Range("B6").Value = "CS304 (Tutorial)"
$wsA.Range("F6").Value = "Free"
$wsA.Range("C7").Value = "Free"
$wsA.Range("E7").Value = "Free"
$wsA.Range("F7").Value = "Free"
$wsA.Range("C8").Value = "CS303 (Tutorial)"
$wsA.Range("E8").Value = "Free"

# --- Section_B sheet updates ---
$wsB = $wb.Worksheets.Item("Section_B")
$wsB.Range("B2").Value = "CS309"
$wsB.Range("E2").Value = "Free"
$wsB.Range("F2").Value = "CS304"
$wsB.Range("B3").Value = "CS303"
$wsB.Range("E3").Value = "CS304"
$wsB.Range("C5").Value = "Free"
$wsB.Range("D5").Value = "Free"
$wsB.Range("F5").Value = "CS309"
$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "CS303 (Tutorial)"
$wsB.Range("F6").Value = "CS304 (Tutorial)"
$wsB.Range("B7").Value = "Free"
$wsB.Range("C7").Value = "CS303"
$wsB.Range("D7").Value = "CS304"
$wsB.Range("E7").Value = "CS303"
$wsB.Range("B8").Value = "Free"
$wsB.Range("D8").Value = "Free"
$wsB.Range("E8").Value = "CS309 (Tutorial)"
